# The TEST_CASES sheet holds a sample row describing a scripted test case.
# Update the test case "kind" cell (X2) from SCRIPTED to GHERKIN so that it
# is consistent with the scripting language / script already recorded in
# the adjoining columns (Y2/Z2), and move the active selection to X3.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TEST_CASES")
$ws.Activate()

$ws.Range("X2").Value = "GHERKIN"

[void]$ws.Range("X3").Select()
